# Figure 6 data update for 3D Notch Test legend workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A/B updates (top block) ---
$ws.Range("A2").Value = 0.2
$ws.Range("A5").Value = 0.05
$ws.Range("A7").Value = -0.44

# --- Column D updates (top block, E has formulas already referencing D) ---
$ws.Range("D10").Value = 1.9
$ws.Range("D11").Value = 1.5
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 0.5

# D14 moves to the small value that used to live in D19, and picks up the
# scientific-notation number format that travelled with it.
$ws.Range("D14").NumberFormat = "0.00E+00"
$ws.Range("D14").Value = -0.0091000000000000004

# --- Column A/B updates (bottom block) ---
$ws.Range("A16").Value = -1.5

# --- Remove the now-unused trailing D/E rows ---
# D15:E15 and D16:E16 become empty (row 15/16 keep their A/B content),
# and rows 17-19 (which only ever held D/E data) disappear entirely.
$ws.Range("D15:E16").ClearContents()
$ws.Range("17:19").Delete()

# --- Restore the selection to the new last-used cell ---
$null = $ws.Range("B19").Select()
